$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.638.79'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').Value = '1.642.46'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.97'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('E8').Value = '  -0.31%  '
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.23'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0842'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('D12').Value = '1.872.69'
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.653.55'
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.21'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.81%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.530'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.57'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.31%  '
$ws.Range('D17').Value = '26.679.35'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '0.0₃0748'
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '216.77'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.36'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.46%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.30'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.91%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.54'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.87%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.17'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +11.13%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '146.03'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.12'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +3.86%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.76'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.25%  '
$ws.Range('E30').Value = '  +2.49%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.17'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('E32').Value = '  +2.41%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.05'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.93%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').Value = '1.268.81'
$ws.Range('E34').Value = '  +4.39%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.53'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.92%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.41'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0180'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +5.11%  '
$ws.Range('E38').Value = '  +5.88%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.826'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.35%  '
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.812'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.43%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.24'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.52%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.46'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.58%  '
$ws.Range('D44').Value = '1.782.59'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '92.92'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.15%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '59.30'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +7.68%  '
$ws.Range('E47').Value = '  +2.51%  '
$ws.Range('E48').Value = '  +1.02%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.79'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.30%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0971'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.91%  '
